$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1723128898356615
$ws.Range("C2").Value = 0.003303860599588913

$ws.Range("B3").Value = 0.3580959746890076
$ws.Range("C3").Value = 0.002839634754745392

$ws.Range("B4").Value = 0.3402074837213397
$ws.Range("C4").Value = 0.004339140753139611

$ws.Range("B5").Value = 0.1624886556157242
$ws.Range("C5").Value = 0.002529849915437787

$ws.Range("B6").Value = 0.1283414277294919
$ws.Range("C6").Value = 0.002374545865387838
